$d = $word.ActiveDocument

$old = "с использованием форматированого вывода."
$new = "с использованием форматированого вывода, а также запустил её на сервере helios."

$d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
